$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels: "2024e" -> "2024" and "2025f" -> "2025e"
$ws.Range("E1").Value = "2024"
$ws.Range("F1").Value = "2025e"
